$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value edits -----------------------------------------------------
# F3: Kinetic energy 20 -> 10 MeV
$ws.Range("F3").Value = 10

# F6 / F7: SigmaX / SigmaY - switch to scientific-notation display format
$ws.Range("F6").NumberFormat = "0.0000E+00"
$ws.Range("F7").NumberFormat = "0.0000E+00"

# F8: MeanEnergy 15 -> 10 MeV
$ws.Range("F8").Value = 10

# F11: Length of first drift 4.034E-2 -> 4.118E-2
$ws.Range("F11").Value = 0.04118

# F18: Gap between colimator first (F)quad and second (D)quad 2.577E-2 -> 3.6953E-2
$ws.Range("F18").Value = 0.036953

# F22: Main drift from last quad to kapton/aluminium foils 1.728652 -> 1.6
$ws.Range("F22").Value = 1.6

# --- Column width ------------------------------------------------------
# Column F widens slightly to fit the new values
$ws.Columns.Item(6).ColumnWidth = 9.5

# --- Selection -----------------------------------------------------------
$null = $ws.Range("F9").Select()
